$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: DXINCRE -> DXINCRE_LEFT ---------------------------------------
$ws.Range("A11").Value = "DXINCRE_LEFT"
$ws.Range("D11").Value = "size increase ratio for the spatial mesh, used for the region to the left of the refined region."

# --- Insert two new rows after row 11 for DXINCRE_RIGHT and MAXNOD --------
$ws.Rows("12:13").Insert(-4121) | Out-Null

# Row 12 picks up row 11's formatting automatically (copy just to be safe)
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A12:E12").PasteSpecial(-4122) | Out-Null

# Row 13 should look like row 7 (plain label/unit/type row, no wrap height)
# for columns A-D, but column E needs the plain "general" style (like column B/C)
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A13:D13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null

# --- Row 12 content: DXINCRE_RIGHT ------------------------------------------
$ws.Range("A12").Value = "DXINCRE_RIGHT"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "float"
$ws.Range("D12").Value = "size increase ratio for the spatial mesh, used for the region to the right of the refined region."
$ws.Range("E12").Value = 1.2

# --- Row 13 content: MAXNOD -------------------------------------------------
$ws.Range("A13").Value = "MAXNOD"
$ws.Range("B13").Value = "-"
$ws.Range("C13").Value = "integer"
$ws.Range("D13").Value = " maximum number of nodes for conductor spatial discretization"
$ws.Range("E13").Value = 10001

$ws.Range("A3").Select() | Out-Null
